$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 0.002
$ws.Range("D2").Value = 0.25
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 250
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0.1
$ws.Range("L2").Value = 0.01
$ws.Range("M2").Value = 15
